# Alumni Locator workbook update: add 3 new members (Chamith, Hansani, Gayashan)
# and correct a couple of existing coordinate values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Fix existing coordinates
# ------------------------------------------------------------------
# Thamali (row 17) - corrected coordinate
$ws.Range("C17").Value2 = -23.2744
$ws.Range("D17").Value2 = 133.7751

# Row 42 - corrected latitude
$ws.Range("C42").Value2 = -21.8688

# ------------------------------------------------------------------
# 2. Add three new rows (46, 47, 48) cloning the look & feel of the
#    most recently added block (row 43/44 style: s=3/4/5/6/7)
# ------------------------------------------------------------------
$ws.Range("A43:G43").Copy() | Out-Null
$ws.Range("A46:G48").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 46 - Chamith
$ws.Range("A46").Value2 = "Chamith"
$ws.Range("B46").Value2 = "Sri Lanka"
$ws.Range("C46").Value2 = 6.9271
$ws.Range("D46").Value2 = 79.8612
$ws.Range("E46").Value2 = "https://www.linkedin.com/in/chamith-nadeeshan/"
$ws.Range("F46").Value2 = "Chamith.jpg"

# Row 47 - Hansani
$ws.Range("A47").Value2 = "Hansani"
$ws.Range("B47").Value2 = "Sri Lanka"
$ws.Range("C47").Value2 = 6.9271
$ws.Range("D47").Value2 = 79.8612
$ws.Range("E47").Value2 = "https://www.linkedin.com/in/hansani-gunathilaka/"
$ws.Range("F47").Value2 = "Hansani .jpg"

# Row 48 - Gayashan
$ws.Range("A48").Value2 = "Gayashan"
$ws.Range("B48").Value2 = "Sri Lanka"
$ws.Range("C48").Value2 = 6.9271
$ws.Range("D48").Value2 = 79.8612
$ws.Range("E48").Value2 = "https://www.linkedin.com/in/pasindu-gayashan/"
$ws.Range("F48").Value2 = "Gayashan.jpg"

# ------------------------------------------------------------------
# 3. Hyperlink the new "Linkedin URL" cells (this resets their style
#    to the built-in Hyperlink style, so we re-stamp the original
#    formatting on top right afterwards).
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E46"), "https://www.linkedin.com/in/chamith-nadeeshan/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E47"), "https://www.linkedin.com/in/hansani-gunathilaka/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E48"), "https://www.linkedin.com/in/pasindu-gayashan/") | Out-Null

$ws.Range("E43").Copy() | Out-Null
$ws.Range("E46:E48").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 4. Update selection / view to rest on the newly added last row
# ------------------------------------------------------------------
$ws.Rows.Item(48).Select() | Out-Null
